# Reorganize the header row:
#   old layout: A1=Email ID, B1=Company Name, C1=Domain
#   new layout: A1=Company Name, B1=First Name, C1=Last Name, D1=Email ID, E1=Domain
#
# First stage the two brand-new headers into the columns that will be
# vacated (D/E) so the shared-string table grows in the same order the
# workbook shows (Email ID, Company Name, Domain, First Name, Last Name),
# then overwrite every header cell with its final value so the old
# strings are dropped from the shared-string table and the new
# A1..E1 values point at the right shared-string indices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "First Name"
$ws.Range("E1").Value = "Last Name"

$ws.Range("A1").Value = "Company Name"
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Last Name"
$ws.Range("D1").Value = "Email ID"
$ws.Range("E1").Value = "Domain"

# Resize the columns to fit their new header text.
$ws.Range("A1:E1").EntireColumn.AutoFit() | Out-Null

# Match the saved selection (active cell D1) recorded in the sheet view.
$ws.Range("D1").Select() | Out-Null
